# fixed UI of all mails - append 16 new "All Mail" rows (S.No 25-40)
# to the "Download - All Mail" sheet, matching the existing table layout
# (A: S.No, B: Receiver Email, C: Sender Name, D: Sender Email, E: Subject).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(26,1).Value = 25
$ws.Cells.Item(26,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(26,3).Value = "Veer_New1"
$ws.Cells.Item(26,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(26,5).Value = "Test-Check-23"
$ws.Cells.Item(27,1).Value = 26
$ws.Cells.Item(27,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(27,3).Value = "Veer_New1"
$ws.Cells.Item(27,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(27,5).Value = "Test-Check-23"
$ws.Cells.Item(28,1).Value = 27
$ws.Cells.Item(28,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(28,3).Value = "Veer_New-01"
$ws.Cells.Item(28,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(28,5).Value = "Test-Check-33"
$ws.Cells.Item(29,1).Value = 28
$ws.Cells.Item(29,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(29,3).Value = "Veer_New-02"
$ws.Cells.Item(29,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(29,5).Value = "Test-Check-33"
$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(30,3).Value = "Veer_New-02"
$ws.Cells.Item(30,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(30,5).Value = "Test-Check-33"
$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(31,3).Value = "Veer_New-02"
$ws.Cells.Item(31,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(31,5).Value = "Test-Check-33"
$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(32,3).Value = "Veer_New-02"
$ws.Cells.Item(32,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(32,5).Value = "Test-Check-33"
$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(33,3).Value = "Veer_New-02"
$ws.Cells.Item(33,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(33,5).Value = "Test-Check-33"
$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(34,3).Value = "Veer_New-03"
$ws.Cells.Item(34,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(34,5).Value = "Test-Check-34"
$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(35,3).Value = "Veer_New-03"
$ws.Cells.Item(35,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(35,5).Value = "Test-Check-34"
$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(36,3).Value = "Veer_New-03"
$ws.Cells.Item(36,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(36,5).Value = "Test-Check-34"
$ws.Cells.Item(37,1).Value = 36
$ws.Cells.Item(37,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(37,3).Value = "Veer_New-03"
$ws.Cells.Item(37,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(37,5).Value = "Test-Check-34"
$ws.Cells.Item(38,1).Value = 37
$ws.Cells.Item(38,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(38,3).Value = "Veer_New-03"
$ws.Cells.Item(38,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(38,5).Value = "Test-Check-36"
$ws.Cells.Item(39,1).Value = 38
$ws.Cells.Item(39,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(39,3).Value = "Veer_New-03"
$ws.Cells.Item(39,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(39,5).Value = "Test-Check-36"
$ws.Cells.Item(40,1).Value = 39
$ws.Cells.Item(40,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(40,3).Value = "Veer_New-03"
$ws.Cells.Item(40,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(40,5).Value = "Test-Check-37"
$ws.Cells.Item(41,1).Value = 40
$ws.Cells.Item(41,2).Value = "veer.edu028@gmail.com"
$ws.Cells.Item(41,3).Value = "Veer_New-03"
$ws.Cells.Item(41,4).Value = "veer.prakash_cs.aiml19@gla.ac.in"
$ws.Cells.Item(41,5).Value = "Test-Check-38"

Write-Host "Appended rows 26-41 to sheet $($ws.Name)"
